# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1097
$ws1.Range("F3").Value = 4098
$ws1.Range("F5").Value = 324
$ws1.Range("F12").Value = 227
$ws1.Range("F13").Value = 2874
$ws1.Range("F15").Value = 1300

# Sheet "全部类型" (all types, contains an extra row shifting the later entries down by one)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1097
$ws4.Range("F3").Value = 4098
$ws4.Range("F5").Value = 324
$ws4.Range("F13").Value = 227
$ws4.Range("F14").Value = 2874
$ws4.Range("F16").Value = 1300
